$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (G=38956)
$ws.Range("H17").Value = 987.45715
$ws.Range("J17").Value = 987.45715
$ws.Range("L17").Value = 2962.37145
$ws.Range("N17").Value = -3298.37145

# Row 28 (G=27772)
$ws.Range("H28").Value = 1708.15
$ws.Range("I28").Value = 1267.6666
$ws.Range("K28").Value = 1267.6666
$ws.Range("M28").Value = -782.6666

# Row 40 (G=5505)
$ws.Range("H40").Value = 19254396
$ws.Range("J40").Value = 41706184
$ws.Range("L40").Value = 41706184
$ws.Range("N40").Value = -41706534

# Row 86 (G=12603)
$ws.Range("H86").Value = 1298464.5
$ws.Range("I86").Value = 2025445
$ws.Range("J86").Value = 6054.778
$ws.Range("K86").Value = 2025445
$ws.Range("L86").Value = 6054.778
$ws.Range("M86").Value = -2024322
$ws.Range("N86").Value = -8300.778

# Row 89 (G=12603)
$ws.Range("H89").Value = 1298464.5
$ws.Range("I89").Value = 2025445
$ws.Range("J89").Value = 6054.778
$ws.Range("K89").Value = 10127225
$ws.Range("L89").Value = 30273.89
$ws.Range("M89").Value = -10121609
$ws.Range("N89").Value = -41505.89

# Row 121 (G=39731)
$ws.Range("H121").Value = 4245.033
$ws.Range("J121").Value = 4245.033
$ws.Range("L121").Value = 12735.099
$ws.Range("N121").Value = -16229.099

# Row 132 (G=44049)
$ws.Range("H132").Value = 14193.206
$ws.Range("I132").Value = 2804.8
$ws.Range("K132").Value = 8414.400000000001
$ws.Range("M132").Value = -5884.400000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (G=27713)
$ws.Range("H2").Value = 1165780.8
$ws.Range("I2").Value = 1456422.2
$ws.Range("K2").Value = 1456422.2
$ws.Range("M2").Value = -1456309.2

# Row 45 (G=27714)
$ws.Range("H45").Value = 43548.57
$ws.Range("I45").Value = 56568.43
$ws.Range("J45").Value = 4489
$ws.Range("K45").Value = 56568.43
$ws.Range("L45").Value = 4489
$ws.Range("M45").Value = -56191.43
$ws.Range("N45").Value = -5243

# Row 74 (G=44000)
$ws.Range("H74").Value = 967.8889
$ws.Range("I74").Value = 742.2353000000001
$ws.Range("J74").Value = 1351.5
$ws.Range("K74").Value = 742.2353000000001
$ws.Range("L74").Value = 1351.5
$ws.Range("M74").Value = 131.7646999999999
$ws.Range("N74").Value = -3099.5

# Row 77 (G=44000)
$ws.Range("H77").Value = 967.8889
$ws.Range("I77").Value = 742.2353000000001
$ws.Range("J77").Value = 1351.5
$ws.Range("K77").Value = 3711.1765
$ws.Range("L77").Value = 6757.5
$ws.Range("M77").Value = 656.8234999999995
$ws.Range("N77").Value = -15493.5

# Row 116 (G=27713)
$ws.Range("H116").Value = 1165780.8
$ws.Range("I116").Value = 1456422.2
$ws.Range("K116").Value = 1456422.2
$ws.Range("M116").Value = -1454128.2

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (G=27713)
$ws.Range("H3").Value = 1165780.8
$ws.Range("I3").Value = 1456422.2
$ws.Range("K3").Value = 1456422.2
$ws.Range("M3").Value = -1456308.2

# Row 93 (G=19510)
$ws.Range("H93").Value = 109975
$ws.Range("J93").Value = 109975
$ws.Range("L93").Value = 109975
$ws.Range("N93").Value = -113719

# Row 97 (G=19518)
$ws.Range("H97").Value = 11245.25
$ws.Range("I97").Value = 11245.25
$ws.Range("K97").Value = 11245.25
$ws.Range("M97").Value = -10254.25

# Row 99 (G=19943)
$ws.Range("H99").Value = 1737568.5
$ws.Range("I99").Value = 2316135.8
$ws.Range("K99").Value = 2316135.8
$ws.Range("M99").Value = -2314637.8

# Row 118 (G=27137)
$ws.Range("H118").Value = 55000
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 55000
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 55000
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -58314

# Row 134 (G=43998)
$ws.Range("H134").Value = 2993.2
$ws.Range("I134").Value = 1226.8235
$ws.Range("K134").Value = 3680.4705
$ws.Range("M134").Value = -1145.4705

$ws = $wb.Worksheets.Item("CRP")
# Row 53 (G=25632)
$ws.Range("H53").Value = 37090.668
$ws.Range("J53").Value = 37090.668
$ws.Range("L53").Value = 37090.668
$ws.Range("N53").Value = -38304.668

# Row 86 (G=12584)
$ws.Range("H86").Value = 9797.75
$ws.Range("I86").Value = 9221
$ws.Range("K86").Value = 9221
$ws.Range("M86").Value = -8098

# Row 89 (G=12584)
$ws.Range("H89").Value = 9797.75
$ws.Range("I89").Value = 9221
$ws.Range("K89").Value = 46105
$ws.Range("M89").Value = -40489

# Row 107 (G=27689)
$ws.Range("H107").Value = 1515690.8
$ws.Range("I107").Value = 3030716.2
$ws.Range("J107").Value = 665.1667
$ws.Range("K107").Value = 3030716.2
$ws.Range("L107").Value = 665.1667
$ws.Range("M107").Value = -3028796.2
$ws.Range("N107").Value = -4505.1667

# Row 134 (G=44020)
$ws.Range("H134").Value = 2449.5
$ws.Range("I134").Value = 1917.5454
$ws.Range("J134").Value = 3619.8
$ws.Range("K134").Value = 5752.6362
$ws.Range("L134").Value = 10859.4
$ws.Range("M134").Value = -3217.6362
$ws.Range("N134").Value = -15929.4

$ws = $wb.Worksheets.Item("CUL")
# Row 11 (G=4745)
$ws.Range("H11").Value = 82081170
$ws.Range("I11").Value = 114521.43
$ws.Range("K11").Value = 343564.29
$ws.Range("M11").Value = -343424.29

# Row 26 (G=4746)
$ws.Range("H26").Value = 163.4
$ws.Range("J26").Value = 502
$ws.Range("L26").Value = 1506
$ws.Range("N26").Value = -2082

# Row 86 (G=12892)
$ws.Range("H86").Value = 300
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

# Row 89 (G=12892)
$ws.Range("H89").Value = 300
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

# Row 104 (G=19807)
$ws.Range("H104").Value = 14318.091
$ws.Range("I104").Value = 2499
$ws.Range("K104").Value = 7497
$ws.Range("M104").Value = -4876

# Row 113 (G=27843)
$ws.Range("H113").Value = 764.375
$ws.Range("I113").Value = 700.625
$ws.Range("J113").Value = 828.125
$ws.Range("K113").Value = 2101.875
$ws.Range("L113").Value = 2484.375
$ws.Range("M113").Value = 68.125
$ws.Range("N113").Value = -6824.375

$ws = $wb.Worksheets.Item("GSM")
# Row 132 (G=44008)
$ws.Range("H132").Value = 160698.92
$ws.Range("I132").Value = 290758.56
$ws.Range("K132").Value = 872275.6799999999
$ws.Range("M132").Value = -869745.6799999999

$ws = $wb.Worksheets.Item("LTW")
# Row 132 (G=44058)
$ws.Range("H132").Value = 3318.7192
$ws.Range("I132").Value = 2512.8223
$ws.Range("K132").Value = 7538.466899999999
$ws.Range("M132").Value = -5008.466899999999

# Row 138 (G=42334)
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280

# Row 141 (G=42487)
$ws.Range("H141").Value = 98500
$ws.Range("J141").Value = 98500
$ws.Range("L141").Value = 98500
$ws.Range("N141").Value = -108860

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (G=12596)
$ws.Range("H81").Value = 10425416
$ws.Range("I81").Value = 10426666
$ws.Range("J81").Value = 10424166
$ws.Range("K81").Value = 20853332
$ws.Range("L81").Value = 20848332
$ws.Range("M81").Value = -20852271
$ws.Range("N81").Value = -20850454

# Row 84 (G=12596)
$ws.Range("H84").Value = 10425416
$ws.Range("I84").Value = 10426666
$ws.Range("J84").Value = 10424166
$ws.Range("K84").Value = 104266660
$ws.Range("L84").Value = 104241660
$ws.Range("M84").Value = -104261356
$ws.Range("N84").Value = -104252268

# Row 132 (G=44029)
$ws.Range("H132").Value = 34729956
$ws.Range("I132").Value = 6174254
$ws.Range("K132").Value = 18522762
$ws.Range("M132").Value = -18520232

Write-Host "Applied all Cactuar_Profits updates"